$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 existing edits: D1 becomes text (leading-zero phone number), F1 gets text format but stays numeric ---
$ws.Range("D1").NumberFormat = "@"
$ws.Range("D1").Value = "0546666666"

$ws.Range("F1").Value = 208063511
$ws.Range("F1").NumberFormat = "@"

# --- Row 2 (B first; C filled in later to reproduce original authoring order) ---
$ws.Range("A2").Value = 112323233
$ws.Range("B2").Value = "av"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0546666666"

$ws.Range("E2").Value = "New"

$ws.Range("F2").Value = 208063511
$ws.Range("F2").NumberFormat = "@"

$ws.Range("G2").Value = 43401
$ws.Range("G2").NumberFormat = "m/d/yy"

$ws.Range("H2").Value = "don’t know what he wants in his life"

# --- Row 3 ---
$ws.Range("A3").Value = 412232131
$ws.Range("B3").Value = "dds"
$ws.Range("C3").Value = "as"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "0546666666"

$ws.Range("E3").Value = "New"

$ws.Range("F3").Value = 208063511
$ws.Range("F3").NumberFormat = "@"

$ws.Range("G3").Value = 43401
$ws.Range("G3").NumberFormat = "m/d/yy"

$ws.Range("H3").Value = "don’t know what he wants in his life"

# --- Row 4 ---
$ws.Range("A4").Value = 123213123
$ws.Range("B4").Value = "dds"
$ws.Range("C4").Value = "as"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0546666666"

$ws.Range("E4").Value = "New"

$ws.Range("F4").Value = 208063511
$ws.Range("F4").NumberFormat = "@"

$ws.Range("G4").Value = 43401
$ws.Range("G4").NumberFormat = "m/d/yy"

$ws.Range("H4").Value = "don’t know what he wants in his life"

# --- Back to Row 2: fill C2 last so "f" lands as the final new shared-string entry ---
$ws.Range("C2").Value = "f"

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection ---
$null = $ws.Range("A6").Select()
